$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantities to standard amounts for 4 servings.
# C2: Eier (Größe M) quantity 2 -> 6
$ws.Range("C2").Value = 6

# C3: Fischsoße quantity "0.5" (text) -> 2 (numeric)
$ws.Range("C3").Value = 2

# C4: Helle Sojasoße quantity "0.5" (text) -> 2 (numeric)
$ws.Range("C4").Value = 2

# Reflect the new active cell/selection (was B2, now C2)
$ws.Range("C2").Select()
